# Update "Name of Algo" result values in column C for specific rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = -12.3846
    14  = -13.922
    16  = -14.35819999999999
    21  = -12.8324
    23  = -12.1377
    25  = -11.5009
    26  = -11.79040000000001
    29  = -12.1707
    40  = -12.0925
    53  = -10.9883
    57  = -14.17609999999999
    59  = -12.58230000000001
    65  = -12.7851
    69  = -11.0649
    79  = -12.2084
    83  = -13.70549999999999
    91  = -12.5983
    93  = -10.3818
    100 = -11.9861
    103 = -12.52479999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
